$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New columns E,F,G,H header row (row 1) - 3 new groups of data: cxq, hyy, hzj
$ws.Range("E1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("F1").Value = "cxq6hz_20170224_144343_ASIC_EEG"
$ws.Range("G1").Value = "hyy-调节6Hz_20170306_110203_ASIC_EEG"
$ws.Range("H1").Value = "hzj-调节6Hz_20170220_113105_ASIC_EEG"

# Row 2 data
$ws.Range("E2").Value = 0.70623145400593468
$ws.Range("F2").Value = 0.70623145400593468
$ws.Range("G2").Value = 0.70707070707070707
$ws.Range("H2").Value = 0.7357357357357357

# Row 3 data
$ws.Range("E3").Value = 0.70845481049562675
$ws.Range("F3").Value = 0.70845481049562675
$ws.Range("G3").Value = 0.72696245733788389
$ws.Range("H3").Value = 0.7129032258064516

# Move selection to match the post-edit cursor position
$ws.Range("I9").Select()

# Best-effort: restore recorded window chrome size (cosmetic, may be a no-op)
$excel.ActiveWindow.Width = 28800
$excel.ActiveWindow.Height = 13170
